$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text for column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Extend column A's date formatting (style from A2) down through A22
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the corrected series of year-end dates (Excel serial numbers) into A2:A22
$dates = @(38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657, 46022)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Write the year-over-year values into B3:B21 (B2 and B22 stay empty)
$values = @{
    3  = 5.361718827437545
    4  = 2.685523658736089
    5  = 0.6019580713040096
    6  = -2.848383350681438
    7  = 0.5670099411379192
    8  = 5.44843673085138
    9  = -0.6152915357131694
    10 = 0.4729742736614195
    11 = 2.566421764830462
    12 = 0.8407878010570302
    13 = 2.194841098049016
    14 = 2.643540836453884
    15 = 2.959935600123309
    16 = 3.452860220335019
    17 = 1.607096457785584
    18 = -3.136134057684858
    19 = 1.613985729693268
    20 = -1.992466799383086
    21 = -3.350381746968589
}
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# B2 previously held the year-over-year value that now belongs to row 3; clear it.
$ws.Range("B2").ClearContents()
